# Applies the "Add files via upload" edit: fills in three previously-empty
# poker-session rows (9, 10, 11 -> worksheet rows 11, 12, 13) on the
# "Spieltabelle" sheet with session data, and restores the view/selection
# state recorded in the saved workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Spieltabelle")

# ---------------------------------------------------------------------
# Row 11 (Index 9)
# ---------------------------------------------------------------------
$ws.Range("B11").Value = "Cashgame"
$ws.Range("C11").Value = "sc.ch"
$ws.Range("D11").Value = 6
$ws.Range("E11").Value = 0.8
$ws.Range("F11").Value = 0.01
$ws.Range("H11").Value = 45955
$ws.Range("I11").Value = 0.70138888888888884
$ws.Range("J11").Value = 45955
$ws.Range("K11").Value = 0.70833333333333337
$ws.Range("M11").Value = 1.5
$ws.Range("N11").Value = 0
$ws.Range("Q11").Value = 1.5
$ws.Range("R11").Value = 3
$ws.Range("S11").Value = 0.66
$ws.Range("T11").Value = "1 Schotch"
$ws.Range("U11").Value = "Hatte ein Paar mit gutem Kicker, Gegner hatte 2 Paare ich ging all in"

# ---------------------------------------------------------------------
# Row 12 (Index 10)
# ---------------------------------------------------------------------
$ws.Range("B12").Value = "Cashgame"
$ws.Range("C12").Value = "sc.ch"
$ws.Range("D12").Value = 6
$ws.Range("E12").Value = 0.8
$ws.Range("F12").Value = 0.01
$ws.Range("H12").Value = 45955
$ws.Range("I12").Value = 0.70833333333333337
$ws.Range("J12").Value = 45955
$ws.Range("K12").Value = 0.73958333333333337
$ws.Range("M12").Value = 1.5
$ws.Range("N12").Value = 0
$ws.Range("Q12").Value = 1.76
$ws.Range("R12").Value = 3
$ws.Range("S12").Value = 0.21
$ws.Range("T12").Value = "1 Schotch"
$ws.Range("U12").Value = "Hatte strasse Gegner hatte flush, konnte aber sehen was passiert"
$ws.Range("V12").Value = "Sehr viel reraisen, bei tight spiel (über 100% des Pots). Bei All in stellte sich heraus, Gegner hatten beim Flop nichts, nur mit Glück gewonnen"

# ---------------------------------------------------------------------
# Row 13 (Index 11)
# ---------------------------------------------------------------------
$ws.Range("B13").Value = "Cashgame"
$ws.Range("C13").Value = "sc.ch"
$ws.Range("D13").Value = 6
$ws.Range("E13").Value = 0.8
$ws.Range("F13").Value = 0.01
$ws.Range("H13").Value = 45955
$ws.Range("I13").Value = 0.73958333333333337
$ws.Range("J13").Value = 45955
$ws.Range("K13").Value = 0.77777777777777779
$ws.Range("M13").Value = 1.5
$ws.Range("N13").Value = 4.92
$ws.Range("Q13").Value = 4.92
$ws.Range("R13").Value = 3
$ws.Range("S13").Value = 0.26
$ws.Range("T13").Value = "2 Schotch"
$ws.Range("U13").Value = "2 Paar Könige siegten gegen 2 tight Player All Inn"
$ws.Range("V13").Value = "Ein Fisch wurde früh entdeckt"

# ---------------------------------------------------------------------
# Restore the sheet's scroll position / selection as last saved.
# ---------------------------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 3
$excel.ActiveWindow.ScrollRow = 1
$excel.Goto($ws.Range("N2:N16"), $false)
